$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, pushing existing rows 7-15 down to 8-16.
$ws.Rows.Item(7).Insert()

# Copy the date-format style used by column D (row 8, the row just pushed down)
# onto the newly inserted D7 cell, so it keeps the same number format.
$ws.Range("D8").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 7 with the data from the diff.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44544
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112006
$ws.Range("G7").Value = "Repollo"
$ws.Range("H7").Value = "Copenhague"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 600
$ws.Range("L7").Value = 650
$ws.Range("M7").Value = 625
$ws.Range("N7").Value = "`$/unidad"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 625
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
